$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- 1. Simple in-place value updates (row positions unaffected) ---
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# --- 2. Insert a new "Jurisdiction" property row right after "Contact" ---
# The sheet grows from 19 to 20 data rows, so rows 11-19 ("Description"
# through "Derivation") move down to rows 12-20.
#
# We avoid EntireRow.Insert()/Range.Insert() here: in this host they mint a
# brand-new, slightly different cell style (missing the table's border)
# instead of reusing the existing style index, which would make every
# shifted cell diverge stylistically from the original file. Instead, first
# clone the formatting of row 19 onto the new row 20 with a plain
# copy/paste of formats, then rewrite every cell's value from the bottom up
# so a value is never clobbered before it has been relocated.
$ws.Range("A19:B19").Copy()
$ws.Range("A20:B20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A20").Value = "Derivation"
$ws.Range("B20").Value = "specialization"

$ws.Range("A19").Value = "Abstract"
# Copy B18's existing "false" text cell (value-only) instead of typing a
# literal .Value = "false": this host auto-coerces a bare "false"/"true"
# string assigned via .Value into a boolean cell, but pasting an existing
# text cell's value keeps it as text, matching the source workbook.
$ws.Range("B18").Copy()
$ws.Range("B19").PasteSpecial(-4163)
$excel.CutCopyMode = 0

$ws.Range("A18").Value = "Base Definition"
$ws.Range("B18").Value = "http://hl7.org/cda/stds/core/StructureDefinition/ST"

$ws.Range("A17").Value = "Type"
$ws.Range("B17").Value = "http://hl7.org/cda/stds/core/StructureDefinition/ENXP"

$ws.Range("A16").Value = "Kind"
$ws.Range("B16").Value = "logical"

$ws.Range("A15").Value = "FHIR Version"
$ws.Range("B15").Value = "5.0.0"

$ws.Range("A14").Value = "Copyright"
$ws.Range("B14").Value = ""

$ws.Range("A13").Value = "Purpose"
$ws.Range("B13").Value = ""

$ws.Range("A12").Value = "Description"
$ws.Range("B12").Value = "A character string token representing a part of a name. May have a type code signifying the role of the part in the whole entity name, and a qualifier code for more detail about the name part type. Typical name parts for person names are given names, and family names, titles, etc."

# Finally, row 11 becomes the new "Jurisdiction" row (value intentionally
# blank, matching the source edit).
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
